# Scheduled runner update: refresh market-price derived columns (H-N) across
# several Leve profit sheets. Values come from the latest Universalis price
# snapshot; LeveProfit columns (M/N) are recomputed accordingly and may be
# added or removed depending on whether a profit/loss applies to a row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H33").Value = 364.84616
$ws.Range("I33").Value = 294.1111
$ws.Range("J33").Value = 524
$ws.Range("K33").Value = 294.1111
$ws.Range("L33").Value = 524
$ws.Range("M33").Value = -65.11110000000002
$ws.Range("N33").Value = -982

$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10748

$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35736

$ws.Range("H86").Value = 5108.5835
$ws.Range("I86").Value = 5479.8
$ws.Range("J86").Value = 4843.4287
$ws.Range("K86").Value = 5479.8
$ws.Range("L86").Value = 4843.4287
$ws.Range("M86").Value = -4356.8
$ws.Range("N86").Value = -7089.4287

$ws.Range("H89").Value = 5108.5835
$ws.Range("I89").Value = 5479.8
$ws.Range("J89").Value = 4843.4287
$ws.Range("K89").Value = 27399
$ws.Range("L89").Value = 24217.1435
$ws.Range("M89").Value = -21783
$ws.Range("N89").Value = -35449.14350000001

$ws.Range("H92").Value = 425.8
$ws.Range("I92").Value = 383.66666
$ws.Range("J92").Value = 805
$ws.Range("K92").Value = 383.66666
$ws.Range("L92").Value = 805
$ws.Range("M92").Value = 864.33334

$ws.Range("H100").Value = 4614.2856
$ws.Range("I100").Value = 5460
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 5460
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -4919

$ws.Range("H106").Value = 1325.75
$ws.Range("I106").Value = 1325.75
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1325.75
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -694.75

$ws.Range("H132").Value = 4670.2915
$ws.Range("I132").Value = 1235.6428
$ws.Range("J132").Value = 9478.799999999999
$ws.Range("K132").Value = 3706.9284
$ws.Range("L132").Value = 28436.4
$ws.Range("M132").Value = -1176.9284
$ws.Range("N132").Value = -33496.39999999999

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 4833.3335
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 6500
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 6500
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = -7254

$ws.Range("H97").Value = 12114
$ws.Range("I97").Value = 153.33333
$ws.Range("J97").Value = 30055
$ws.Range("K97").Value = 153.33333
$ws.Range("L97").Value = 30055
$ws.Range("M97").Value = 342.66667
$ws.Range("N97").Value = -31047

$ws.Range("H131").Value = 90000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 90000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 90000
$ws.Range("N131").Value = -100080

# ---------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws.Range("H94").Value = 704.7778
$ws.Range("I94").Value = 733.8823
$ws.Range("J94").Value = 210
$ws.Range("K94").Value = 733.8823
$ws.Range("L94").Value = 210
$ws.Range("M94").Value = -282.8823
$ws.Range("N94").Value = -1112

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H58").Value = 5372.25
$ws.Range("J58").Value = 11784.5
$ws.Range("L58").Value = 11784.5
$ws.Range("N58").Value = -12190.5

$ws.Range("H59").Value = 43748
$ws.Range("I59").Value = 39999.5
$ws.Range("J59").Value = 47496.5
$ws.Range("K59").Value = 39999.5
$ws.Range("L59").Value = 47496.5
$ws.Range("M59").Value = -38854.5
$ws.Range("N59").Value = -49786.5

$ws.Range("H60").Value = 20516.25
$ws.Range("I60").Value = 4093
$ws.Range("J60").Value = 25990.666
$ws.Range("K60").Value = 4093
$ws.Range("L60").Value = 25990.666
$ws.Range("M60").Value = -3582
$ws.Range("N60").Value = -27012.666

$ws.Range("H103").Value = 5250
$ws.Range("I103").Value = 5250
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 5250
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -4078
$ws.Range("N103").ClearContents()

$ws.Range("H132").Value = 7288.5938
$ws.Range("I132").Value = 5509.077
$ws.Range("K132").Value = 16527.231
$ws.Range("M132").Value = -13997.231

$ws.Range("H136").Value = 5372.25
$ws.Range("J136").Value = 11784.5
$ws.Range("L136").Value = 35353.5
$ws.Range("N136").Value = -40453.5

# ---------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H97").Value = 1125.5555
$ws.Range("I97").Value = 1168.7142
$ws.Range("K97").Value = 1168.7142
$ws.Range("M97").Value = -672.7141999999999

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 1605
$ws.Range("I46").Value = 1166.6666
$ws.Range("J46").Value = 1914.4117
$ws.Range("K46").Value = 1166.6666
$ws.Range("L46").Value = 1914.4117
$ws.Range("M46").Value = -978.6666
$ws.Range("N46").Value = -2290.4117

$ws.Range("H82").Value = 3134.625
$ws.Range("I82").Value = 797.5
$ws.Range("J82").Value = 3913.6667
$ws.Range("K82").Value = 797.5
$ws.Range("L82").Value = 3913.6667
$ws.Range("M82").Value = -436.5
$ws.Range("N82").Value = -4635.6667

$ws.Range("H85").Value = 3134.625
$ws.Range("I85").Value = 797.5
$ws.Range("J85").Value = 3913.6667
$ws.Range("K85").Value = 797.5
$ws.Range("L85").Value = 3913.6667
$ws.Range("M85").Value = 450.5
$ws.Range("N85").Value = -6409.6667

$ws.Range("H93").Value = 1488.0667
$ws.Range("I93").Value = 1491.4445
$ws.Range("K93").Value = 1491.4445
$ws.Range("M93").Value = -243.4445000000001

$ws.Range("H100").Value = 2843.3845
$ws.Range("I100").Value = 1495.8334
$ws.Range("J100").Value = 3998.4285
$ws.Range("K100").Value = 1495.8334
$ws.Range("L100").Value = 3998.4285
$ws.Range("M100").Value = -954.8334
$ws.Range("N100").Value = -5080.4285

# ---------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 850
$ws.Range("J81").Value = 875
$ws.Range("L81").Value = 1750
$ws.Range("N81").Value = -3872

$ws.Range("H84").Value = 850
$ws.Range("J84").Value = 875
$ws.Range("L84").Value = 8750
$ws.Range("N84").Value = -19358

$ws.Range("H96").Value = 2226.182
$ws.Range("J96").Value = 3224.75
$ws.Range("L96").Value = 3224.75
$ws.Range("N96").Value = -5970.75

$ws.Range("H136").Value = 2172.5
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
